# The three observation rows (4, 5 and 6) on the "Artfynd" sheet got their
# data rotated: row 4's data moved to row 6, row 5's data moved to row 4,
# and row 6's data moved to row 5 (i.e. a shift up by one row, wrapping
# around). Only columns A, P, Q, R, S, AW and AX actually hold values that
# differ between the three rows - everything else (B..H, T..AB, AD, AE, AG,
# AT, AY, etc.) is identical across rows 4-6 already.
#
# In addition a handful of otherwise-empty "marker" cells (I, J, K, L, N,
# AF) are present/absent differently per row; those follow the very same
# rotation (row4 <- row5's set, row5 <- row6's set, row6 <- row4's set).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$valueCols = @("A", "P", "Q", "R", "S", "AW", "AX")

# Snapshot "before" values for rows 4, 5 and 6 first, so subsequent writes
# do not clobber data that still needs to be read.
$row4 = @{}
$row5 = @{}
$row6 = @{}
foreach ($col in $valueCols) {
    $row4[$col] = $ws.Range($col + "4").Value2
    $row5[$col] = $ws.Range($col + "5").Value2
    $row6[$col] = $ws.Range($col + "6").Value2
}

# Row 4 becomes old row 5's data.
foreach ($col in $valueCols) {
    $ws.Range($col + "4").Value = $row5[$col]
}

# Row 5 becomes old row 6's data.
foreach ($col in $valueCols) {
    $ws.Range($col + "5").Value = $row6[$col]
}

# Row 6 becomes old row 4's data.
foreach ($col in $valueCols) {
    $ws.Range($col + "6").Value = $row4[$col]
}

# Row 6 previously only had the empty marker cells I6 and K6 (same set as
# row 4 had before the edit); J6, L6, N6 and AF6 need to disappear now.
$ws.Range("J6").Clear()
$ws.Range("L6").Clear()
$ws.Range("N6").Clear()
$ws.Range("AF6").Clear()

# Row 4 now needs the fuller set of empty marker cells that row 5 had
# (I4, J4, K4, L4, N4, AF4). I4 and K4 already existed; add the rest as
# empty cells without introducing any new formatting/style.
$ws.Range("J4").Font.Bold = $false
$ws.Range("L4").Font.Bold = $false
$ws.Range("N4").Font.Bold = $false
$ws.Range("AF4").Font.Bold = $false
